$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(3, 4).Value = 0.18181818181818182
$ws.Cells.Item(3, 5).Value = 2.0
$ws.Cells.Item(4, 6).Value = 0.16666666666666666
$ws.Cells.Item(4, 7).Value = 2.0
$ws.Cells.Item(4, 8).Value = 0.16666666666666666
$ws.Cells.Item(4, 9).Value = 2.0
$ws.Cells.Item(4, 10).Value = 0.08333333333333333
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.08333333333333333
$ws.Cells.Item(4, 13).Value = 1.0
$ws.Cells.Item(4, 14).Value = 0.16666666666666666
$ws.Cells.Item(4, 15).Value = 2.0
$ws.Cells.Item(5, 6).Value = 0.023809523809523808
$ws.Cells.Item(5, 7).Value = 1.0
$ws.Cells.Item(5, 8).Value = 0.07142857142857142
$ws.Cells.Item(5, 9).Value = 3.0
$ws.Cells.Item(5, 14).Value = 0.07142857142857142
$ws.Cells.Item(5, 15).Value = 3.0
$ws.Cells.Item(6, 4).Value = 0.07407407407407407
$ws.Cells.Item(6, 5).Value = 2.0
$ws.Cells.Item(6, 6).Value = 0.07407407407407407
$ws.Cells.Item(6, 7).Value = 2.0
$ws.Cells.Item(6, 8).Value = 0.07407407407407407
$ws.Cells.Item(6, 9).Value = 2.0
$ws.Cells.Item(6, 14).Value = 0.07407407407407407
$ws.Cells.Item(6, 15).Value = 2.0
$ws.Cells.Item(7, 6).Value = 0.07692307692307693
$ws.Cells.Item(7, 7).Value = 1.0
$ws.Cells.Item(7, 8).Value = 0.15384615384615385
$ws.Cells.Item(7, 9).Value = 2.0
$ws.Cells.Item(8, 6).Value = 0.03773584905660377
$ws.Cells.Item(8, 7).Value = 2.0
$ws.Cells.Item(8, 8).Value = 0.05660377358490566
$ws.Cells.Item(8, 9).Value = 3.0
$ws.Cells.Item(8, 10).Value = 0.018867924528301886
$ws.Cells.Item(8, 11).Value = 1.0
$ws.Cells.Item(8, 12).Value = 0.018867924528301886
$ws.Cells.Item(8, 13).Value = 1.0
$ws.Cells.Item(8, 14).Value = 0.03773584905660377
$ws.Cells.Item(8, 15).Value = 2.0
$ws.Cells.Item(9, 6).Value = 0.03333333333333333
$ws.Cells.Item(9, 7).Value = 1.0
$ws.Cells.Item(9, 8).Value = 0.06666666666666667
$ws.Cells.Item(9, 9).Value = 2.0
$ws.Cells.Item(9, 10).Value = 0.03333333333333333
$ws.Cells.Item(9, 11).Value = 1.0
$ws.Cells.Item(10, 4).Value = 0.02564102564102564
$ws.Cells.Item(10, 5).Value = 1.0
$ws.Cells.Item(10, 10).Value = 0.02564102564102564
$ws.Cells.Item(10, 11).Value = 1.0
$ws.Cells.Item(10, 12).Value = 0.07692307692307693
$ws.Cells.Item(10, 13).Value = 3.0
$ws.Cells.Item(10, 14).Value = 0.10256410256410256
$ws.Cells.Item(10, 15).Value = 4.0
$ws.Cells.Item(11, 4).Value = 0.06666666666666667
$ws.Cells.Item(11, 5).Value = 1.0
$ws.Cells.Item(11, 6).Value = 0.06666666666666667
$ws.Cells.Item(11, 7).Value = 1.0
$ws.Cells.Item(11, 8).Value = 0.06666666666666667
$ws.Cells.Item(11, 9).Value = 1.0
$ws.Cells.Item(11, 10).Value = 0.13333333333333333
$ws.Cells.Item(11, 11).Value = 2.0
$ws.Cells.Item(11, 12).Value = 0.13333333333333333
$ws.Cells.Item(11, 13).Value = 2.0
$ws.Cells.Item(11, 14).Value = 0.13333333333333333
$ws.Cells.Item(11, 15).Value = 2.0
$ws.Cells.Item(12, 10).Value = 0.16666666666666666
$ws.Cells.Item(12, 11).Value = 1.0
$ws.Cells.Item(12, 12).Value = 0.16666666666666666
$ws.Cells.Item(12, 13).Value = 1.0
$ws.Cells.Item(12, 14).Value = 0.16666666666666666
$ws.Cells.Item(12, 15).Value = 1.0
$ws.Cells.Item(13, 4).Value = 0.047619047619047616
$ws.Cells.Item(13, 5).Value = 1.0
$ws.Cells.Item(13, 6).Value = 0.09523809523809523
$ws.Cells.Item(13, 7).Value = 2.0
$ws.Cells.Item(13, 8).Value = 0.09523809523809523
$ws.Cells.Item(13, 9).Value = 2.0
$ws.Cells.Item(13, 12).Value = 0.047619047619047616
$ws.Cells.Item(13, 13).Value = 1.0
$ws.Cells.Item(13, 14).Value = 0.09523809523809523
$ws.Cells.Item(13, 15).Value = 2.0
$ws.Cells.Item(14, 6).Value = 0.043478260869565216
$ws.Cells.Item(14, 7).Value = 1.0
$ws.Cells.Item(14, 8).Value = 0.13043478260869565
$ws.Cells.Item(14, 9).Value = 3.0
$ws.Cells.Item(14, 10).Value = 0.043478260869565216
$ws.Cells.Item(14, 11).Value = 1.0
$ws.Cells.Item(14, 12).Value = 0.043478260869565216
$ws.Cells.Item(14, 13).Value = 1.0
$ws.Cells.Item(14, 14).Value = 0.17391304347826086
$ws.Cells.Item(14, 15).Value = 4.0
$ws.Cells.Item(15, 4).Value = 0.03225806451612903
$ws.Cells.Item(15, 5).Value = 1.0
$ws.Cells.Item(15, 6).Value = 0.06451612903225806
$ws.Cells.Item(15, 7).Value = 2.0
$ws.Cells.Item(15, 8).Value = 0.12903225806451613
$ws.Cells.Item(15, 9).Value = 4.0
$ws.Cells.Item(15, 12).Value = 0.06451612903225806
$ws.Cells.Item(15, 13).Value = 2.0
$ws.Cells.Item(15, 14).Value = 0.0967741935483871
$ws.Cells.Item(15, 15).Value = 3.0
$ws.Cells.Item(16, 4).Value = 0.03225806451612903
$ws.Cells.Item(16, 5).Value = 1.0
$ws.Cells.Item(16, 6).Value = 0.06451612903225806
$ws.Cells.Item(16, 7).Value = 2.0
$ws.Cells.Item(16, 10).Value = 0.03225806451612903
$ws.Cells.Item(16, 11).Value = 1.0
$ws.Cells.Item(16, 12).Value = 0.03225806451612903
$ws.Cells.Item(16, 13).Value = 1.0
$ws.Cells.Item(16, 14).Value = 0.06451612903225806
$ws.Cells.Item(16, 15).Value = 2.0
$ws.Cells.Item(19, 4).Value = 0.05263157894736842
$ws.Cells.Item(19, 5).Value = 1.0
$ws.Cells.Item(19, 6).Value = 0.10526315789473684
$ws.Cells.Item(19, 7).Value = 2.0
$ws.Cells.Item(19, 8).Value = 0.21052631578947367
$ws.Cells.Item(19, 9).Value = 4.0
$ws.Cells.Item(19, 10).Value = 0.10526315789473684
$ws.Cells.Item(19, 11).Value = 2.0
$ws.Cells.Item(20, 4).Value = 0.14285714285714285
$ws.Cells.Item(20, 5).Value = 4.0
$ws.Cells.Item(20, 12).Value = 0.03571428571428571
$ws.Cells.Item(20, 13).Value = 1.0
$ws.Cells.Item(20, 14).Value = 0.07142857142857142
$ws.Cells.Item(20, 15).Value = 2.0
$ws.Cells.Item(22, 12).Value = 0.05263157894736842
$ws.Cells.Item(22, 13).Value = 1.0
$ws.Cells.Item(24, 6).Value = 0.125
$ws.Cells.Item(24, 7).Value = 1.0
$ws.Cells.Item(24, 8).Value = 0.25
$ws.Cells.Item(24, 9).Value = 2.0
$ws.Cells.Item(29, 4).Value = 0.10344827586206896
$ws.Cells.Item(29, 5).Value = 3.0
$ws.Cells.Item(29, 10).Value = 0.034482758620689655
$ws.Cells.Item(29, 11).Value = 1.0
$ws.Cells.Item(29, 12).Value = 0.10344827586206896
$ws.Cells.Item(29, 13).Value = 3.0
$ws.Cells.Item(32, 4).Value = 0.03571428571428571
$ws.Cells.Item(32, 5).Value = 2.0
$ws.Cells.Item(32, 10).Value = 0.03571428571428571
$ws.Cells.Item(32, 11).Value = 2.0
$ws.Cells.Item(32, 12).Value = 0.05357142857142857
$ws.Cells.Item(32, 13).Value = 3.0
$ws.Cells.Item(35, 8).Value = 0.2
$ws.Cells.Item(35, 9).Value = 1.0
$ws.Cells.Item(35, 14).Value = 0.2
$ws.Cells.Item(35, 15).Value = 1.0
$ws.Cells.Item(37, 14).Value = 0.1111111111111111
$ws.Cells.Item(37, 15).Value = 2.0
$ws.Cells.Item(38, 4).Value = 0.07142857142857142
$ws.Cells.Item(38, 5).Value = 1.0
$ws.Cells.Item(38, 6).Value = 0.07142857142857142
$ws.Cells.Item(38, 7).Value = 1.0
$ws.Cells.Item(38, 8).Value = 0.07142857142857142
$ws.Cells.Item(38, 9).Value = 1.0
$ws.Cells.Item(41, 4).Value = 0.08571428571428572
$ws.Cells.Item(41, 5).Value = 3.0
$ws.Cells.Item(41, 6).Value = 0.11428571428571428
$ws.Cells.Item(41, 7).Value = 4.0
$ws.Cells.Item(41, 10).Value = 0.02857142857142857
$ws.Cells.Item(41, 11).Value = 1.0
$ws.Cells.Item(41, 12).Value = 0.05714285714285714
$ws.Cells.Item(41, 13).Value = 2.0
$ws.Cells.Item(41, 14).Value = 0.11428571428571428
$ws.Cells.Item(41, 15).Value = 4.0
$ws.Cells.Item(43, 6).Value = 0.03125
$ws.Cells.Item(43, 7).Value = 1.0
$ws.Cells.Item(43, 8).Value = 0.125
$ws.Cells.Item(43, 9).Value = 4.0
$ws.Cells.Item(43, 10).Value = 0.03125
$ws.Cells.Item(43, 11).Value = 1.0
$ws.Cells.Item(43, 12).Value = 0.03125
$ws.Cells.Item(43, 13).Value = 1.0
$ws.Cells.Item(43, 14).Value = 0.125
$ws.Cells.Item(43, 15).Value = 4.0
$ws.Cells.Item(44, 4).Value = 0.125
$ws.Cells.Item(44, 5).Value = 3.0
$ws.Cells.Item(44, 6).Value = 0.16666666666666666
$ws.Cells.Item(44, 7).Value = 4.0
$ws.Cells.Item(44, 10).Value = 0.041666666666666664
$ws.Cells.Item(44, 11).Value = 1.0
$ws.Cells.Item(44, 12).Value = 0.08333333333333333
$ws.Cells.Item(44, 13).Value = 2.0
$ws.Cells.Item(44, 14).Value = 0.125
$ws.Cells.Item(44, 15).Value = 3.0
$ws.Cells.Item(45, 4).Value = 0.12
$ws.Cells.Item(45, 5).Value = 3.0
$ws.Cells.Item(45, 6).Value = 0.16
$ws.Cells.Item(45, 7).Value = 4.0
$ws.Cells.Item(45, 10).Value = 0.04
$ws.Cells.Item(45, 11).Value = 1.0
$ws.Cells.Item(45, 12).Value = 0.08
$ws.Cells.Item(45, 13).Value = 2.0
$ws.Cells.Item(46, 4).Value = 0.047619047619047616
$ws.Cells.Item(46, 5).Value = 1.0
$ws.Cells.Item(46, 6).Value = 0.09523809523809523
$ws.Cells.Item(46, 7).Value = 2.0
$ws.Cells.Item(46, 8).Value = 0.14285714285714285
$ws.Cells.Item(46, 9).Value = 3.0
$ws.Cells.Item(46, 10).Value = 0.09523809523809523
$ws.Cells.Item(46, 11).Value = 2.0
$ws.Cells.Item(46, 12).Value = 0.09523809523809523
$ws.Cells.Item(46, 13).Value = 2.0
$ws.Cells.Item(46, 14).Value = 0.19047619047619047
$ws.Cells.Item(46, 15).Value = 4.0
$ws.Cells.Item(47, 4).Value = 0.027777777777777776
$ws.Cells.Item(47, 5).Value = 1.0
$ws.Cells.Item(47, 6).Value = 0.05555555555555555
$ws.Cells.Item(47, 7).Value = 2.0
$ws.Cells.Item(47, 8).Value = 0.08333333333333333
$ws.Cells.Item(47, 9).Value = 3.0
$ws.Cells.Item(47, 10).Value = 0.027777777777777776
$ws.Cells.Item(47, 11).Value = 1.0
$ws.Cells.Item(47, 12).Value = 0.05555555555555555
$ws.Cells.Item(47, 13).Value = 2.0
$ws.Cells.Item(47, 14).Value = 0.05555555555555555
$ws.Cells.Item(47, 15).Value = 2.0
$ws.Cells.Item(48, 4).Value = 0.0625
$ws.Cells.Item(48, 5).Value = 3.0
$ws.Cells.Item(48, 6).Value = 0.0625
$ws.Cells.Item(48, 7).Value = 3.0
$ws.Cells.Item(49, 4).Value = 0.07407407407407407
$ws.Cells.Item(49, 5).Value = 4.0
$ws.Cells.Item(49, 10).Value = 0.037037037037037035
$ws.Cells.Item(49, 11).Value = 2.0
$ws.Cells.Item(49, 12).Value = 0.037037037037037035
$ws.Cells.Item(49, 13).Value = 2.0
$ws.Cells.Item(49, 14).Value = 0.07407407407407407
$ws.Cells.Item(49, 15).Value = 4.0
$ws.Cells.Item(50, 4).Value = 0.07317073170731707
$ws.Cells.Item(50, 5).Value = 3.0
$ws.Cells.Item(50, 6).Value = 0.0975609756097561
$ws.Cells.Item(50, 7).Value = 4.0
$ws.Cells.Item(50, 10).Value = 0.04878048780487805
$ws.Cells.Item(50, 11).Value = 2.0
$ws.Cells.Item(51, 4).Value = 0.03225806451612903
$ws.Cells.Item(51, 5).Value = 1.0
$ws.Cells.Item(51, 6).Value = 0.0967741935483871
$ws.Cells.Item(51, 7).Value = 3.0
$ws.Cells.Item(51, 10).Value = 0.06451612903225806
$ws.Cells.Item(51, 11).Value = 2.0
$ws.Cells.Item(51, 12).Value = 0.0967741935483871
$ws.Cells.Item(51, 13).Value = 3.0
$ws.Cells.Item(51, 14).Value = 0.12903225806451613
$ws.Cells.Item(51, 15).Value = 4.0
$ws.Cells.Item(52, 4).Value = 0.11764705882352941
$ws.Cells.Item(52, 5).Value = 2.0
$ws.Cells.Item(52, 6).Value = 0.17647058823529413
$ws.Cells.Item(52, 7).Value = 3.0
$ws.Cells.Item(52, 8).Value = 0.17647058823529413
$ws.Cells.Item(52, 9).Value = 3.0
$ws.Cells.Item(52, 12).Value = 0.058823529411764705
$ws.Cells.Item(52, 13).Value = 1.0
$ws.Cells.Item(52, 14).Value = 0.058823529411764705
$ws.Cells.Item(52, 15).Value = 1.0
$ws.Cells.Item(54, 6).Value = 0.03225806451612903
$ws.Cells.Item(54, 7).Value = 1.0
$ws.Cells.Item(54, 8).Value = 0.0967741935483871
$ws.Cells.Item(54, 9).Value = 3.0
$ws.Cells.Item(54, 12).Value = 0.06451612903225806
$ws.Cells.Item(54, 13).Value = 2.0
